$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'36.500.93"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.954.05"
$ws.Range("E3").Value = "  +0.85%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +2.92%  "

# Row 7 - Solana
$ws.Range("D7").Value = "'60.08"
$ws.Range("E7").Value = "  +6.62%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +5.17%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.0787"
$ws.Range("E10").Value = "  -2.33%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.02%  "

# Row 12 - Chainlink
$ws.Range("D12").Value = "'14.16"
$ws.Range("E12").Value = "  +7.30%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "'0.840"
$ws.Range("E13").Value = "  +5.17%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'2.241.98"
$ws.Range("E14").Value = "  +0.88%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "'21.54"
$ws.Range("E15").Value = "  +0.87%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "'5.26"
$ws.Range("E16").Value = "  +2.73%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'1.960.04"
$ws.Range("E17").Value = "  +1.54%  "

# Row 18
$ws.Range("D18").Value = "'36.442.93"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19
$ws.Range("D19").Value = "'69.22"
$ws.Range("E19").Value = "  +0.50%  "

# Row 20
$ws.Range("E20").Value = "  +0.41%  "

# Row 21
$ws.Range("D21").Value = "'229.26"
$ws.Range("E21").Value = "  +1.29%  "

# Row 22
$ws.Range("E22").Value = "  +3.24%  "

# Row 23
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("E24").Value = "  +3.18%  "

# Row 25
$ws.Range("E25").Value = "  +3.97%  "

# Row 26
$ws.Range("E26").Value = "  +7.44%  "

# Row 27
$ws.Range("E27").Value = "  +1.18%  "

# Row 28
$ws.Range("D28").Value = "'160.46"
$ws.Range("E28").Value = "  +0.72%  "

# Row 29
$ws.Range("D29").Value = "'19.27"
$ws.Range("E29").Value = "  +1.28%  "

# Row 30
$ws.Range("E30").Value = "  +21.50%  "

# Row 31
$ws.Range("E31").Value = "  +2.13%  "

# Row 32
$ws.Range("E32").Value = "  +4.66%  "

# Row 33
$ws.Range("E33").Value = "  +0.45%  "

# Row 34
$ws.Range("E34").Value = "  +8.42%  "

# Row 35
$ws.Range("E35").Value = "  +0.05%  "

# Row 36 - becomes LidoDAOToken (was RenderToken)
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.26"
$ws.Range("E36").Value = "  +3.46%  "

# Row 37 - becomes RenderToken (was LidoDAOToken)
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'3.41"
$ws.Range("E37").Value = "  +6.90%  "

# Row 38
$ws.Range("E38").Value = "  -0.44%  "

# Row 39
$ws.Range("D39").Value = "'5.44"
$ws.Range("E39").Value = "  -10.75%  "

# Row 40
$ws.Range("D40").Value = "'0.0965"
$ws.Range("E40").Value = "  -1.42%  "

# Row 41
$ws.Range("E41").Value = "  +0.38%  "

# Row 42
$ws.Range("E42").Value = "  +2.60%  "

# Row 43
$ws.Range("E43").Value = "  +1.14%  "

# Row 44
$ws.Range("D44").Value = "'15.84"
$ws.Range("E44").Value = "  +1.26%  "

# Row 45
$ws.Range("D45").Value = "'1.360.48"
$ws.Range("E45").Value = "  +2.41%  "

# Row 46
$ws.Range("D46").Value = "'88.75"
$ws.Range("E46").Value = "  +3.94%  "

# Row 47
$ws.Range("D47").Value = "'1.02"
$ws.Range("E47").Value = "  +0.63%  "

# Row 48
$ws.Range("E48").Value = "  +1.22%  "

# Row 49
$ws.Range("E49").Value = "  +0.83%  "

# Row 50
$ws.Range("D50").Value = "'46.06"
$ws.Range("E50").Value = "  +7.65%  "

# Row 51
$ws.Range("D51").Value = "'2.137.14"
$ws.Range("E51").Value = "  +1.06%  "
